# Weekly data refresh: a new week's price record is inserted at the top of
# the data table (row 15, just below the fixed first block of rows 2-14),
# every existing data row shifts down by one, and the oldest row spills into
# a brand-new row 47 at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 15
$lastDataRow  = 46
$newLastRow   = 47
$lastCol      = 18   # column R

# Shift rows [firstDataRow .. lastDataRow] down into [firstDataRow+1 .. newLastRow],
# working from the bottom up so a row isn't overwritten before it's copied.
# Column D (the date column) carries an explicit date NumberFormat that must
# travel with the value; every other column uses the sheet's default format,
# so leave those cells alone to avoid minting needless new styles.
$dateCol = 4
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell  = $ws.Cells.Item($r, $c)
        $destCell = $ws.Cells.Item($destRow, $c)
        if ($c -eq $dateCol) {
            $destCell.NumberFormat = $srcCell.NumberFormat
        }
        $destCell.Value = $srcCell.Value()
    }
}

# Populate the new top row (row 15) with this week's record.
$ws.Cells.Item($firstDataRow, 4).Value  = 44544   # D: Fecha
$ws.Cells.Item($firstDataRow, 9).Value  = "Primera" # I: Calidad (unchanged)
$ws.Cells.Item($firstDataRow, 10).Value = 250       # J: Volumen (unchanged)
$ws.Cells.Item($firstDataRow, 11).Value = 1000      # K: Precio minimo
$ws.Cells.Item($firstDataRow, 12).Value = 1200      # L: Precio maximo
$ws.Cells.Item($firstDataRow, 13).Value = 1100      # M: Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 16).Value = 367        # P: Precio $/Kg
